$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C; this shifts the old
# DESIGN CODE..CATALOGS columns (C..G) one place to the right (D..H)
# and keeps their widths/values intact.
$ws.Columns("C").Insert()

# New "ALT CODE" column header and column width (matches neighboring
# OUR CODE / DESIGN CODE columns). Excel's ColumnWidth (chars) stores as
# width + ~0.8333 (the default-font padding), so subtract that constant
# to land on an exact stored width.
$padding = 0.8333333333333333
$ws.Range("C1").Value = "ALT CODE"
$ws.Columns("C").ColumnWidth = 15 - $padding

# Former GSM/CATALOGS columns (now F/G after the insert) change width.
$ws.Columns("F").ColumnWidth = 15 - $padding
$ws.Columns("G").ColumnWidth = 10 - $padding

# Fill in the new ALT CODE value for the existing product row, and
# update its CATALOGS value.
$ws.Range("C2").Value = "501W"
$ws.Range("H2").Value = "Woodrica"

# Duplicate row 2 into a new row 3 (copies values + formatting, so
# text-like values such as "60" keep being stored as text instead of
# being reinterpreted as numbers), then patch the two cells that
# differ for the new "501A" / "Artvio" variant.
$ws.Range("A2:H2").Copy($ws.Range("A3"))
$ws.Range("B3").Value = "501A"
$ws.Range("H3").Value = "Artvio"
